$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.08%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.76%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.069"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.99%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07813"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.19%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.260"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.17%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.087"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.79%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.046"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.90%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9294"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.68%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09418"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.52%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1829"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.94%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08548"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.36%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.73%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09954"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.85%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001476"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.14%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005714"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.02%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.98%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-5.34%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.24%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1321"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.17%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.555"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.27%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.04%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04666"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.85%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.23%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004545"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.09%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.62%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-20.02%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.42%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04710"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.08%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007949"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.00%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.01%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007995"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-18.23%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002292"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.59%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009074"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.97%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006192"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.81%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.61%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.069"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "59.51%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002691"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.27%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.61%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.61%"
